# Apply profit/price recalculations to several Leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 497.625
$ws.Range("I12").Value = 489
$ws.Range("J12").Value = 506.25
$ws.Range("K12").Value = 489
$ws.Range("L12").Value = 506.25
$ws.Range("M12").Value = -319
$ws.Range("N12").Value = -846.25
$ws.Range("H95").Value = 26892.25
$ws.Range("J95").Value = 26892.25
$ws.Range("L95").Value = 26892.25
$ws.Range("N95").Value = -32384.25
$ws.Range("H96").Value = 6254913
$ws.Range("I96").Value = 6679.375
$ws.Range("K96").Value = 20038.125
$ws.Range("M96").Value = -18665.125
$ws.Range("H116").Value = 7436.143
$ws.Range("I116").Value = 7018
$ws.Range("K116").Value = 7018
$ws.Range("M116").Value = -3576
$ws.Range("H137").Value = 8633350
$ws.Range("I137").Value = 25002024
$ws.Range("J137").Value = 18257.684
$ws.Range("K137").Value = 75006072
$ws.Range("L137").Value = 54773.052
$ws.Range("M137").Value = -75003522
$ws.Range("N137").Value = -59873.052
$ws.Range("H138").Value = 2581.55
$ws.Range("I138").Value = 1762.1
$ws.Range("J138").Value = 2854.7
$ws.Range("K138").Value = 5286.299999999999
$ws.Range("L138").Value = 8564.099999999999
$ws.Range("M138").Value = -146.2999999999993
$ws.Range("N138").Value = -18844.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 2505500
$ws.Range("I6").Value = 2505500
$ws.Range("K6").Value = 2505500
$ws.Range("M6").Value = -2505327
$ws.Range("H94").Value = 40665.6
$ws.Range("J94").Value = 40665.6
$ws.Range("L94").Value = 40665.6
$ws.Range("N94").Value = -42467.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 83627.2
$ws.Range("I86").Value = 3755.25
$ws.Range("J86").Value = 225621.78
$ws.Range("K86").Value = 3755.25
$ws.Range("L86").Value = 225621.78
$ws.Range("M86").Value = -2632.25
$ws.Range("N86").Value = -227867.78
$ws.Range("H89").Value = 83627.2
$ws.Range("I89").Value = 3755.25
$ws.Range("J89").Value = 225621.78
$ws.Range("K89").Value = 18776.25
$ws.Range("L89").Value = 1128108.9
$ws.Range("M89").Value = -13160.25
$ws.Range("N89").Value = -1139340.9
$ws.Range("H95").Value = 104372.6
$ws.Range("J95").Value = 104372.6
$ws.Range("L95").Value = 104372.6
$ws.Range("N95").Value = -109864.6
$ws.Range("H107").Value = 2923
$ws.Range("I107").Value = 1558.64
$ws.Range("K107").Value = 1558.64
$ws.Range("M107").Value = 361.3599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 608.1667
$ws.Range("J7").Value = 816.3333
$ws.Range("L7").Value = 816.3333
$ws.Range("N7").Value = -1042.3333
$ws.Range("H14").Value = 505
$ws.Range("I14").Value = 500
$ws.Range("J14").Value = 510
$ws.Range("K14").Value = 500
$ws.Range("L14").Value = 510
$ws.Range("M14").Value = -330
$ws.Range("N14").Value = -850
$ws.Range("H16").Value = 4333.75
$ws.Range("I16").Value = 3741.0908
$ws.Range("J16").Value = 5058.1113
$ws.Range("K16").Value = 3741.0908
$ws.Range("L16").Value = 5058.1113
$ws.Range("M16").Value = -3454.0908
$ws.Range("N16").Value = -5632.1113
$ws.Range("H31").Value = 52637268
$ws.Range("I31").Value = 100002560
$ws.Range("K31").Value = 100002560
$ws.Range("M31").Value = -100002265
$ws.Range("H34").Value = 52637268
$ws.Range("I34").Value = 100002560
$ws.Range("K34").Value = 100002560
$ws.Range("M34").Value = -100002358
$ws.Range("H43").Value = 77780.78
$ws.Range("J43").Value = 77780.78
$ws.Range("L43").Value = 77780.78
$ws.Range("N43").Value = -78148.78
$ws.Range("H50").Value = 38998.332
$ws.Range("J50").Value = 69995
$ws.Range("L50").Value = 69995
$ws.Range("N50").Value = -71245
$ws.Range("H59").Value = 1000000000
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H60").Value = 125023520
$ws.Range("J60").Value = 333366530
$ws.Range("L60").Value = 333366530
$ws.Range("N60").Value = -333367552
$ws.Range("H101").Value = 77780.78
$ws.Range("J101").Value = 77780.78
$ws.Range("L101").Value = 77780.78
$ws.Range("N101").Value = -84270.78
$ws.Range("H105").Value = 1661.6
$ws.Range("I105").Value = 1661.6
$ws.Range("K105").Value = 1661.6
$ws.Range("M105").Value = 85.40000000000009
$ws.Range("H113").Value = 4333.75
$ws.Range("I113").Value = 3741.0908
$ws.Range("J113").Value = 5058.1113
$ws.Range("K113").Value = 3741.0908
$ws.Range("L113").Value = 5058.1113
$ws.Range("M113").Value = -1571.0908
$ws.Range("N113").Value = -9398.1113
$ws.Range("H122").Value = 56499.5
$ws.Range("I122").Value = 63449.5
$ws.Range("J122").Value = 899.5
$ws.Range("K122").Value = 190348.5
$ws.Range("L122").Value = 2698.5
$ws.Range("M122").Value = -187898.5
$ws.Range("N122").Value = -7598.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1045
$ws.Range("J107").Value = 1233.1818
$ws.Range("L107").Value = 3699.5454
$ws.Range("N107").Value = -7539.5454
$ws.Range("H140").Value = 1194.4615
$ws.Range("I140").Value = 1003.5
$ws.Range("J140").Value = 1500
$ws.Range("K140").Value = 3010.5
$ws.Range("L140").Value = 4500
$ws.Range("M140").Value = 2169.5
$ws.Range("N140").Value = -14860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I122").Value = 8320.333000000001
$ws.Range("J122").Value = 5999
$ws.Range("K122").Value = 24960.999
$ws.Range("L122").Value = 17997
$ws.Range("M122").Value = -22510.999
$ws.Range("N122").Value = -22897

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12628.786
$ws.Range("I7").Value = 9980.299999999999
$ws.Range("K7").Value = 9980.299999999999
$ws.Range("M7").Value = -9868.299999999999
$ws.Range("H46").Value = 2635905.2
$ws.Range("I46").Value = 9092197
$ws.Range("J46").Value = 5564.148
$ws.Range("K46").Value = 9092197
$ws.Range("L46").Value = 5564.148
$ws.Range("M46").Value = -9092009
$ws.Range("N46").Value = -5940.148
$ws.Range("H93").Value = 2312.9375
$ws.Range("I93").Value = 1477.1666
$ws.Range("J93").Value = 3387.5
$ws.Range("K93").Value = 1477.1666
$ws.Range("L93").Value = 3387.5
$ws.Range("M93").Value = -229.1666
$ws.Range("N93").Value = -5883.5
$ws.Range("H126").Value = 12628.786
$ws.Range("I126").Value = 9980.299999999999
$ws.Range("K126").Value = 29940.9
$ws.Range("M126").Value = -27470.9
